$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: Model Accuracy (-0.35, 0.35, 0.35)
# Add new columns C:G (Market threshold, Market min, Market max,
# Recall, Precision), update the Accuracy (%) column B values, and
# fill in the new metric columns for every ticker row.
# ---------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Model Accuracy (-0.35, 0.35, 0.35)")

# New header cells, styled like the existing "Accuracy (%)" header (B1).
$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)

$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

# Row 2: TOTALENERGIES SE
$wsAcc.Range("B2").Value = 58.80195599022004
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 0
$wsAcc.Range("G2").Value = 0

# Row 3: FMC CORP
$wsAcc.Range("B3").Value = 33.67970660146699
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 6.970509383378016
$wsAcc.Range("G3").Value = 25

# Row 4: BP PLC
$wsAcc.Range("B4").Value = 87.89731051344744
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

# Row 5: STORA ENSO
$wsAcc.Range("B5").Value = 75.91687041564792
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 1.818181818181818
$wsAcc.Range("G5").Value = 6.451612903225806

# Row 6: BHP GROUP
$wsAcc.Range("B6").Value = 88.93643031784842
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 0
$wsAcc.Range("G6").Value = 0

# ---------------------------------------------------------------
# Sheet 2: Confusion Matrix TOTALENERGIES SE
# Only the "Predicted Neutral" row changes.
# ---------------------------------------------------------------
$wsCm1 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.35, 0.35, 0.35)")
$wsCm1.Range("B3").Value = 9
$wsCm1.Range("C3").Value = 958
$wsCm1.Range("D3").Value = 6

# ---------------------------------------------------------------
# Sheet 3: Confusion Matrix FMC CORP
# ---------------------------------------------------------------
$wsCm2 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.35, 0.35, 0.35)")
$wsCm2.Range("B2").Value = 26
$wsCm2.Range("C2").Value = 54
$wsCm2.Range("D2").Value = 24

$wsCm2.Range("B3").Value = 234
$wsCm2.Range("C3").Value = 428
$wsCm2.Range("D3").Value = 235

$wsCm2.Range("B4").Value = 113
$wsCm2.Range("C4").Value = 170
$wsCm2.Range("D4").Value = 97

# ---------------------------------------------------------------
# Sheet 4: Confusion Matrix BP PLC
# ---------------------------------------------------------------
$wsCm3 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.35, 0.35, 0.35)")
$wsCm3.Range("B3").Value = 38
$wsCm3.Range("C3").Value = 1435
$wsCm3.Range("D3").Value = 39

$wsCm3.Range("B4").Value = 2
$wsCm3.Range("C4").Value = 81

# ---------------------------------------------------------------
# Sheet 5: Confusion Matrix STORA ENSO
# ---------------------------------------------------------------
$wsCm4 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.35, 0.35, 0.35)")
$wsCm4.Range("B2").Value = 2
$wsCm4.Range("C2").Value = 28

$wsCm4.Range("B3").Value = 100
$wsCm4.Range("C3").Value = 1236
$wsCm4.Range("D3").Value = 102

$wsCm4.Range("B4").Value = 8
$wsCm4.Range("C4").Value = 94
$wsCm4.Range("D4").Value = 4

# ---------------------------------------------------------------
# Sheet 6: Confusion Matrix BHP GROUP
# ---------------------------------------------------------------
$wsCm5 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.35, 0.35, 0.35)")
$wsCm5.Range("B2").Value = 0
$wsCm5.Range("C2").Value = 72

$wsCm5.Range("B3").Value = 4
$wsCm5.Range("C3").Value = 1455
$wsCm5.Range("D3").Value = 2

$wsCm5.Range("B4").Value = 0
$wsCm5.Range("C4").Value = 46
